# Update symbol list values (Price and Volume(1h)) per commit
# "Updated symbol list on Sat Feb 11 21:43:29 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.41"
$ws.Range("E2").Value = "'0.83%"
$ws.Range("D3").Value = "'40.95"
$ws.Range("E3").Value = "'0.99%"
$ws.Range("D4").Value = "'5.121"
$ws.Range("E4").Value = "'1.37%"
$ws.Range("D5").Value = "'0.07627"
$ws.Range("E5").Value = "'0.62%"
$ws.Range("D6").Value = "'4.284"
$ws.Range("E6").Value = "'0.35%"
$ws.Range("D7").Value = "'1.604"
$ws.Range("E7").Value = "'0.40%"
$ws.Range("D9").Value = "'0.9101"
$ws.Range("E9").Value = "'0.78%"
$ws.Range("D10").Value = "'0.1261"
$ws.Range("E10").Value = "'23.13%"
$ws.Range("E11").Value = "'2.77%"
$ws.Range("D12").Value = "'0.09095"
$ws.Range("E12").Value = "'0.34%"
$ws.Range("D13").Value = "'0.04340"
$ws.Range("E13").Value = "'2.37%"
$ws.Range("E14").Value = "'-0.64%"
$ws.Range("D15").Value = "'0.001247"
$ws.Range("E15").Value = "'0.56%"
$ws.Range("D16").Value = "'0.005809"
$ws.Range("E16").Value = "'-0.98%"
$ws.Range("E17").Value = "'-0.04%"
$ws.Range("E18").Value = "'1.53%"
$ws.Range("D19").Value = "'6.939"
$ws.Range("E19").Value = "'2.49%"
$ws.Range("E20").Value = "'2.42%"
$ws.Range("D21").Value = "'0.2738"
$ws.Range("E21").Value = "'0.18%"
$ws.Range("E22").Value = "'-3.33%"
$ws.Range("E23").Value = "'3.41%"
$ws.Range("D24").Value = "'0.004103"
$ws.Range("E24").Value = "'1.10%"
$ws.Range("E25").Value = "'-2.55%"
$ws.Range("E26").Value = "'24.24%"
$ws.Range("D38").Value = "'0.02423"
$ws.Range("E38").Value = "'1.89%"
$ws.Range("D39").Value = "'0.05240"
$ws.Range("E39").Value = "'1.59%"
$ws.Range("D40").Value = "'0.007834"
$ws.Range("E40").Value = "'0.75%"
$ws.Range("E41").Value = "'1.16%"
$ws.Range("D42").Value = "'0.006811"
$ws.Range("E42").Value = "'-3.79%"
$ws.Range("E43").Value = "'-3.25%"
$ws.Range("D44").Value = "'0.007452"
$ws.Range("E44").Value = "'-12.45%"
$ws.Range("E45").Value = "'0.07%"
$ws.Range("D46").Value = "'0.00006874"
$ws.Range("E46").Value = "'7.92%"
$ws.Range("E47").Value = "'-0.28%"
$ws.Range("D48").Value = "'0.1442"
$ws.Range("E48").Value = "'2,119.01%"
$ws.Range("E49").Value = "'-32.03%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.28%"
$ws.Range("E51").Value = "'-0.28%"
